# Sixth Session Day 6
#
# On the "Installation" slide, the paragraph that used to read:
#   "...Sometimes if initialization not happened ngOnInit method..."
# is edited so it now reads:
#   "...If initialization not happened ngOnInit method..."
#
# i.e. the words "Sometimes if " (right after the second manual line
# break) are replaced with "If ".

$p = $ppt.ActivePresentation

$needle = "Sometimes if initialization not happened"
$targetSlide = $null
$targetShape = $null

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text.Contains($needle)) {
                    $targetSlide = $s
                    $targetShape = $shp
                }
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not locate the shape containing '$needle'"
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$startPos = $fullText.IndexOf("Sometimes if ") + 1
$oldLen = "Sometimes if ".Length

$rng = $tr.Characters($startPos, $oldLen)
$rng.Text = "If "
